$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.14243932011947891
$ws.Range("A2").Value = -0.097974783404570331
$ws.Range("A3").Value = -0.0089999990776465921
$ws.Range("A4").Value = 0.28398730965901464
$ws.Range("A5").Value = -0.0059999990998225172
$ws.Range("A6").Value = -0.005999999065572581
$ws.Range("A7").Value = -0.019999998918052597
$ws.Range("A8").Value = -0.019999998909956851
$ws.Range("A9").Value = -0.0059999990458603492
$ws.Range("A10").Value = -0.0059999990387638036
$ws.Range("A11").Value = -0.0044999990534009271
$ws.Range("A12").Value = -0.0059999990366876865
$ws.Range("A13").Value = -0.0059999990314771878
$ws.Range("A14").Value = -0.011999998968402537
$ws.Range("A15").Value = -0.0059999990303847284
$ws.Range("A16").Value = 0.022603005201248028
$ws.Range("A17").Value = -0.0059999990276766724
$ws.Range("A18").Value = -0.0089999989959768101
$ws.Range("A19").Value = -0.0089999990787581474
$ws.Range("A20").Value = -0.0089999990705944555
$ws.Range("A21").Value = -0.0089999990695304177
$ws.Range("A22").Value = -0.0089999990687985587
$ws.Range("A23").Value = -0.0068480847665002287
$ws.Range("A24").Value = -0.041999998713937003
$ws.Range("A25").Value = -0.041999998706962138
$ws.Range("A26").Value = -0.0059999990621371069
$ws.Range("A27").Value = -0.0059999990575372308
$ws.Range("A28").Value = -0.018012105345378338
$ws.Range("A29").Value = -0.011999998964004277
$ws.Range("A30").Value = -0.019999998876166991
$ws.Range("A31").Value = -0.014999998921069846
$ws.Range("A32").Value = -0.012335243385160233
$ws.Range("A33").Value = -0.0059999990130119585
